# edit.ps1 - applies the LOB1258.docx content changes via Word COM-interop
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute(
        $find,      # FindText
        $true,      # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $replace,   # ReplaceWith
        2           # Replace (wdReplaceAll)
    ) | Out-Null
}

Replace-Text "Ativação: 01/01/2022" "Ativação: 01/01/2024"

Replace-Text "Fornecer ao aluno os conceitos básicos de hidráulica aplicadas ao meio ambiente." "Fornecer ao aluno os conceitos básicos sobre escoamento de líquidos em conduto forçado"

Replace-Text "Provide the students the basic concepts of hydraulic applied to the environment." "Provide the student with the basic concepts of liquid flow in penstocks."

Replace-Text "Escoamento permanente uniforme em condutos forçados e perdas de carga nas tubulações. Instalações de bombeamento e bombas hidráulicas. Escoamento permanente uniforme em condutos livres, resistência ao escoamento e perdas de carga nos canais naturais e artificiais. Estudo da carga específica em canais." "Propriedades de líquidos. Reologia de fluidos. Escoamento permanente uniforme em condutos forçados de fluidos Newtonianos e não-Newtonianos. Medidores de pressão e vazão. Perdas de carga em tubulações e acessórios. Dimensionamento de bombas e sistemas de bombeamento em instalações industriais."

Replace-Text "Permanent flow in pressure conduits. Head losses in pipelines. Pumping stations and hydraulic pumps. Permanente flow in free surface conduits, head losses in natural and artificial channels. Specific Energy in channels." "Liquid properties. Fluid rheology. Uniform permanent flow in penstocks of Newtonian and non-Newtonian fluids. Pressure and flow gauges. Load losses in pipes and fittings. Sizing of pumps and pumping system in industrial installations."

Replace-Text "- Hidrostática,- piezometria,- conservação da massa e quantidade de movimento,- Escoamentos Permanentes em Condutos Forçados,- Resistência ao Escoamento e Perdas de Carga,- Bombas e sistemas de recalque,- Escoamento Permanente Uniforme em Condutos Livres,- Resistência ao escoamento e Perdas de Carga,- Canais regulares e naturais,- Carga Específica,- Escoamento Permanente Gradualmente Variado,- Cálculo da linha d’água,- Ressalto Hidráulico." "- Propriedades dos líquidos- Reologia de fluidos- Escoamento em conduto forçado de fluidos Newtonianos e não-Newtonianos- Medidores de pressão e vazão- Perdas de carga em tubulações e acessórios para fluidos Newtonianos e não-Newtonianos- Seleção de bombas, classificação e altura de projeto- Curvas características de bombas, cálculo de potência, cálculo do NPSH e cavitação"

Replace-Text "- Hydrostatic,- piezometry,- conservation of mass and amount of movement,- Permanent flows in Pipes,- Resistance to Flow and Load Losses,- Pumps and discharge systems,- Uniform Permanent Flow in Free Flues,- Resistance to flow and head losses,- Regular and natural channels,- Specific Load,- Gradually Variable Permanent Flow,- Calculation of the water line,- Hydraulic boss." "- Properties of liquids- Fluid rheology- Penstock flow of Newtonian and non-Newtonian fluids- Pressure and flow gauges- Pressure drops in pipes and accessories for Newtonian and non-Newtonian fluids- Pump selection, rating, and design height- Pump characteristic curves, power calculation, NPSH calculation and cavitation"

Replace-Text "Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos.Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios." "A avaliação será composta por atividades avaliativas, às quais poderão incluir provas teóricas e trabalhos em grupos, com a elaboração e resolução de problemas reais das instalações industriais."

Replace-Text "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios." "A média do semestre será MF = (M1+M2)/2, sendo que M1 e M2 correspondem às médias parciais obtidas durante o curso, e podem conter provas teóricas e trabalhos, conforme descrito no método acima."

Replace-Text "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0." "1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0"

Replace-Text "Hauer, B.F. Lamberti, G.A. – Methods in Stream Ecology, 896 p, Academic Press, ISBN 0123329078, 2006Alfredini, Paolo; ""Obras e Gestão de Portos e Costas - A Técnica Aliada ao Enfoque Logístico e Ambiental"".Editora Edgard Blucher, São Paulo, 2005.ALFREDINI, P. - “Obras e Gestão de Portos e Costas - A Técnica Aliada ao Enfoque Logístico e Ambiental”. Editora Edgard Blucher, São Paulo, 2005.LENCASTRE, A. – “Hidráulica Geral”. Edição do Autor, Lisboa, 1996.OPEN UNIVERSITY COURSE TEAM – “Waves, tides and shallow water processes”. 1998ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 1999. (Coordenação: Luiz Di Bernardo). ASSOCIAÇÃO BRASILEIRA DE ENGENHARIA SANITÁRIA E AMBIENTAL, Rio de Janeiro, 200 (Coordenação: Marco A.P. Reali)." "MACINTYRE, ARCHIBALD JOSEPH. Manual de Instalações Hidráulicas e Sanitárias eBook Kindle, ed. LTC; 2ª edição, 2020. AZEVEDO NETO, J. M .; FERNANDEZ, M. F.; ITO, A. E., Manual de Hidráulica. 8ª ed. SãoPaulo: Edgar Blucher, 2000.BAPTISTA, M. B.; COELHO, M. M. L. P.; Fundamentos de Engenharia Hidráulica. 3ª ed. Belo Horizonte: Editora UFMG, 2010."
